$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text & layout -------------------------------------------------
# Update the confidence-interval caption in the merged header row.
$ws.Range("A1").Value = "Based on 64 simulations and with 95.00% confidence"

# The merged header used to span A1:E1; it now only spans A1:D1, freeing up
# E1. Give E1 a "vertical center only" alignment (no horizontal centering).
$ws.Range("A1:E1").UnMerge()
$ws.Range("A1:D1").Merge()
$ws.Range("E1").HorizontalAlignment = 1
$ws.Range("E1").VerticalAlignment = -4108

# Update the sheet selection to the new header range.
$ws.Range("A1:D1").Select()

# --- Updated simulation results (64 runs / 95% confidence) ---------------
$ws.Range("B4").Value = 0.52643300000000004
$ws.Range("D4").Value = 0.00033599999999999998
$ws.Range("B5").Value = 2.7959870000000002
$ws.Range("D5").Value = 0.0073680000000000004
$ws.Range("B6").Value = 0.79675499999999999
$ws.Range("D6").Value = 0.0066470000000000001
$ws.Range("B7").Value = 1.9992319999999999
$ws.Range("D7").Value = 0.001248
$ws.Range("B8").Value = 5.310511
$ws.Range("D8").Value = 0.015446
$ws.Range("B9").Value = 1.5133369999999999
$ws.Range("D9").Value = 0.013006
$ws.Range("B10").Value = 0.75943499999999997
$ws.Range("D10").Value = 0.000669

$ws.Range("B14").Value = 0.84750999999999999
$ws.Range("D14").Value = 0.00073399999999999995
$ws.Range("B15").Value = 3.540006
$ws.Range("D15").Value = 0.005058
$ws.Range("B16").Value = 0.338115
$ws.Range("D16").Value = 0.0034710000000000001
$ws.Range("B17").Value = 3.2018909999999998
$ws.Range("D17").Value = 0.0026919999999999999
$ws.Range("B18").Value = 4.1765400000000001
$ws.Range("D18").Value = 0.0079719999999999999
$ws.Range("B19").Value = 0.39893400000000001
$ws.Range("D19").Value = 0.0042570000000000004
$ws.Range("B20").Value = 0.62960099999999997
$ws.Range("D20").Value = 0.00078700000000000005

$ws.Range("B24").Value = 1.321704
$ws.Range("D24").Value = 0.0014760000000000001
$ws.Range("B25").Value = 3.3829750000000001
$ws.Range("D25").Value = 0.0098329999999999997
$ws.Range("B26").Value = 0.88562600000000002
$ws.Range("D26").Value = 0.0078639999999999995
$ws.Range("B27").Value = 2.4973489999999998
$ws.Range("D27").Value = 0.0026970000000000002
$ws.Range("B28").Value = 2.5594459999999999
$ws.Range("D28").Value = 0.0090100000000000006
$ws.Range("B29").Value = 0.670072
$ws.Range("D29").Value = 0.0063449999999999999
$ws.Range("B30").Value = 0.62979099999999999
$ws.Range("D30").Value = 0.001023

$ws.Range("B34").Value = 0.38653999999999999
$ws.Range("D34").Value = 0.000195
$ws.Range("B35").Value = 2.1085820000000002
$ws.Range("D35").Value = 0.0031519999999999999
$ws.Range("B36").Value = 0.80889800000000001
$ws.Range("D36").Value = 0.002516
$ws.Range("B37").Value = 1.2996840000000001
$ws.Range("D37").Value = 0.00079699999999999997
$ws.Range("B38").Value = 5.4549209999999997
$ws.Range("D38").Value = 0.0086119999999999999
$ws.Range("B39").Value = 2.092632
$ws.Range("D39").Value = 0.0067149999999999996
$ws.Range("B40").Value = 0.84057199999999999
$ws.Range("D40").Value = 0.00053700000000000004
$ws.Range("B41").Value = 0.041416000000000001
$ws.Range("D41").Value = 0.00025599999999999999

$ws.Range("B43").Value = 11.827500000000001
$ws.Range("D43").Value = 0.0149
